# Add handlers for adding employee, product
$wb = $excel.ActiveWorkbook

# All cells in this workbook are stored as plain text (Telegram-bot export
# convention), so numeric-looking values are entered with a leading
# apostrophe to force text storage instead of being auto-detected as numbers.

# Sheet: Cотрудники (Employees) - add two new rows (13, 14)
$wsEmployees = $wb.Worksheets.Item("Cотрудники")
$wsEmployees.Range("A13").Value = "'87654321"
$wsEmployees.Range("B13").Value = "Ысаков Акылбек"
$wsEmployees.Range("A14").Value = "'123"
$wsEmployees.Range("B14").Value = "Пол полыч"

# Sheet: Товары (Products) - add a new row (2)
$wsProducts = $wb.Worksheets.Item("Товары")
$wsProducts.Range("A2").Value = "'12345"
$wsProducts.Range("B2").Value = "Контроллер"
$wsProducts.Range("C2").Value = "'12"
$wsProducts.Range("D2").Value = "Акыл"
$wsProducts.Range("E2").Value = "-"
